# Auto-generated edit script: updates cryptos list (prices & volume %) per commit
# "Updated cryptos list on Tue Oct  8 19:36:30 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain number (e.g. '5.20') need to be forced to
# text first, otherwise Excel auto-converts them to a numeric value and drops
# formatting such as trailing zeros (matches the source sheet, which stores every
# Price/Volume cell as text).
$textForceCells = $ws.Range("D4","D5","D6","D8","D12","D13","D14","D15","D19","D20","D21","D22","D23","D25","D26","D27","D30","D31","D32","D35","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
$textForceCells.NumberFormat = "@"

# --- Apply the new values (grouped by row, in sheet order) ---
# Row 2
$ws.Range("D2").Value = '62.203.80'
$ws.Range("E2").Value = '  -1.72%  '

# Row 3
$ws.Range("D3").Value = '2.442.61'
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = '581.84'
$ws.Range("E5").Value = '  +1.77%  '

# Row 6
$ws.Range("D6").Value = '143.28'
$ws.Range("E6").Value = '  -2.39%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").Value = '  -1.82%  '

# Row 9
$ws.Range("D9").Value = '2.440.52'
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("E10").Value = '  -3.14%  '

# Row 11
$ws.Range("E11").Value = '  +1.95%  '

# Row 12
$ws.Range("D12").Value = '5.20'
$ws.Range("E12").Value = '  -0.64%  '

# Row 13
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").Value = '  -3.15%  '

# Row 14
$ws.Range("D14").Value = '26.39'
$ws.Range("E14").Value = '  -2.32%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  -3.63%  '

# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.878.14'
$ws.Range("E16").Value = '  +0.06%  '

# Row 17
$ws.Range("D17").Value = '62.066.56'
$ws.Range("E17").Value = '  -1.56%  '

# Row 18
$ws.Range("D18").Value = '2.439.22'
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$ws.Range("D19").Value = '10.89'
$ws.Range("E19").Value = '  -3.73%  '

# Row 20
$ws.Range("D20").Value = '7.10'
$ws.Range("E20").Value = '  -3.12%  '

# Row 21
$ws.Range("D21").Value = '330.47'
$ws.Range("E21").Value = '  +0.97%  '

# Row 22
$ws.Range("D22").Value = '4.11'
$ws.Range("E22").Value = '  -1.89%  '

# Row 23
$ws.Range("D23").Value = '1.96'
$ws.Range("E23").Value = '  -6.34%  '

# Row 24
$ws.Range("E24").Value = '  +0.00%  '

# Row 25
$ws.Range("D25").Value = '65.64'
$ws.Range("E25").Value = '  +0.52%  '

# Row 26
$ws.Range("D26").Value = '9.35'
$ws.Range("E26").Value = '  +4.58%  '

# Row 27
$ws.Range("D27").Value = '618.67'
$ws.Range("E27").Value = '  -0.01%  '

# Row 28
$ws.Range("D28").Value = '2.557.22'
$ws.Range("E28").Value = '  -0.16%  '

# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0952'
$ws.Range("E29").Value = '  -7.71%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("D31").Value = '1.43'
$ws.Range("E31").Value = '  -4.43%  '

# Row 32
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").Value = '  -3.60%  '

# Row 33
$ws.Range("E33").Value = '  -0.02%  '

# Row 34
$ws.Range("E34").Value = '  -1.00%  '

# Row 35
$ws.Range("D35").Value = '4.91'
$ws.Range("E35").Value = '  -5.84%  '

# Row 36
$ws.Range("E36").Value = '  +0.16%  '

# Row 37
$ws.Range("D37").Value = '1.43'
$ws.Range("E37").Value = '  -6.33%  '

# Row 38
$ws.Range("D38").Value = '0.376'
$ws.Range("E38").Value = '  -1.13%  '

# Row 39
$ws.Range("D39").Value = '151.19'
$ws.Range("E39").Value = '  +3.41%  '

# Row 40
$ws.Range("D40").Value = '18.33'
$ws.Range("E40").Value = '  -2.18%  '

# Row 41
$ws.Range("D41").Value = '5.23'
$ws.Range("E41").Value = '  -3.46%  '

# Row 42
$ws.Range("D42").Value = '1.76'
$ws.Range("E42").Value = '  -1.71%  '

# Row 43
$ws.Range("D43").Value = '42.48'
$ws.Range("E43").Value = '  +1.50%  '

# Row 44
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").Value = '  -9.06%  '

# Row 46
$ws.Range("D46").Value = '143.10'
$ws.Range("E46").Value = '  -3.60%  '

# Row 47
$ws.Range("D47").Value = '3.62'
$ws.Range("E47").Value = '  -3.55%  '

# Row 48
$ws.Range("D48").Value = '0.0524'
$ws.Range("E48").Value = '  -1.84%  '

# Row 49
$ws.Range("D49").Value = '0.598'
$ws.Range("E49").Value = '  -0.48%  '

# Row 50
$ws.Range("D50").Value = '19.50'
$ws.Range("E50").Value = '  -7.94%  '

# Row 51
$ws.Range("D51").Value = '0.0906'
$ws.Range("E51").Value = '  -1.41%  '

# --- Restore default (unstyled) cell style on the text-forced cells so only the
# cell *content* differs from the original workbook, matching the source diff. ---
$textForceCells.Style = "Normal"

